# Applies the "28 февраля 2024 (среда)" schedule update:
#  - trims stray trailing blank lines from 3 shared text cells
#  - adjusts several row heights (content re-wrapped / re-measured)
#  - updates the saved sheet view (zoom level + selected cell, no frozen scroll position)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Text fixes: strip trailing blank line left over in these cells ---
$ws.Range("B1").Value = "28 февраля 2024 (среда) "
$ws.Range("C10").Value = "Информационные технологии в профессиональной деятельности`nСоколова А.С."
$ws.Range("G20").Value = "Основы анализа бухгалтерской отчетности`nГадомская Т.М."
$ws.Range("G21").Value = "Основы анализа бухгалтерской отчетности`nГадомская Т.М."

# --- Row height adjustments ---
$ws.Rows.Item(5).RowHeight = 81.75
$ws.Rows.Item(6).RowHeight = 63.75
$ws.Rows.Item(7).RowHeight = 63.75
$ws.Rows.Item(8).RowHeight = 64.5
$ws.Rows.Item(10).RowHeight = 98.25
$ws.Rows.Item(11).RowHeight = 58.5
$ws.Rows.Item(12).RowHeight = 81
$ws.Rows.Item(13).RowHeight = 66
$ws.Rows.Item(14).RowHeight = 97.5
$ws.Rows.Item(17).RowHeight = 81
$ws.Rows.Item(18).RowHeight = 82.5
$ws.Rows.Item(19).RowHeight = 69
$ws.Rows.Item(20).RowHeight = 81
$ws.Rows.Item(22).RowHeight = 84
$ws.Rows.Item(24).RowHeight = 32.25
$ws.Rows.Item(28).RowHeight = 78.75

# --- Sheet view: select L3 and zoom to 73% (removes the old scrolled-down view) ---
$ws.Range("L3").Select()
$excel.ActiveWindow.Zoom = 73
